$wb = $excel.ActiveWorkbook

# Sheet order: 1=2025, 2=2030, 3=2035, 4=2040, 5=2045, 6=2050
$sheetData = @(
    @{
        Index = 1
        Cells = @{
            "B2" = 1037.265132737054
            "E2" = 28926.05393052954
            "G2" = 8095.925712661834
            "I2" = 16171.06685703679
            "L2" = 48492.22142001599
            "M2" = 10595.37713982
            "N2" = 7070.228218264912
            "O2" = 6993.478371476022
        }
    },
    @{
        Index = 2
        Cells = @{
            "A2" = 0
            "B2" = 4157.588990853394
            "E2" = 45991.90904307188
            "G2" = 8095.925712661834
            "I2" = 37079.12819938764
            "L2" = 54844.03303316472
            "M2" = 17449.04999683176
            "N2" = 9022.5600255922
            "O2" = 9723.120256756203
        }
    },
    @{
        Index = 3
        Cells = @{
            "A2" = 2754.31755456332
            "B2" = 6368.910634126893
            "E2" = 57457.45307013817
            "G2" = 8095.925712661834
            "I2" = 52465.73681402855
            "L2" = 54844.03303316472
            "M2" = 21912.87293902603
            "N2" = 13031.04653826744
            "O2" = 12858.90853139509
        }
    },
    @{
        Index = 4
        Cells = @{
            "A2" = 2754.31755456332
            "B2" = 6368.910634126893
            "E2" = 57457.45307013817
            "G2" = 8095.925712661834
            "I2" = 52465.73681402855
            "L2" = 54844.03303316472
            "M2" = 21912.87293902603
            "N2" = 13148.52234297407
            "O2" = 12858.90853139509
        }
    },
    @{
        Index = 5
        Cells = @{
            "A2" = 5713.151062849596
            "B2" = 6368.910634126893
            "E2" = 57457.45307013817
            "G2" = 8095.925712661834
            "I2" = 52465.73681402855
            "L2" = 54844.03303316472
            "M2" = 21912.87293902603
            "N2" = 13597.36669969182
            "O2" = 14934.99104992778
        }
    },
    @{
        Index = 6
        Cells = @{
            "A2" = 5713.151062849596
            "B2" = 6368.910634126893
            "E2" = 57457.45307013817
            "G2" = 8095.925712661834
            "I2" = 52465.73681402855
            "L2" = 54844.03303316472
            "M2" = 21912.87293902603
            "N2" = 13597.36669969182
            "O2" = 14934.99104992778
        }
    }
)

foreach ($sheetEntry in $sheetData) {
    $ws = $wb.Worksheets.Item($sheetEntry.Index)
    $cellValues = $sheetEntry.Cells
    foreach ($cellRef in $cellValues.Keys) {
        $ws.Range($cellRef).Value = $cellValues[$cellRef]
    }
}
